$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update data values per the diff
$ws.Range("B4").Value = 12
$ws.Range("F5").Value = 6
$ws.Range("F7").Value = 15
$ws.Range("B8").Value = 15

# Update selection
$ws.Range("G11").Select()
